$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look like plain numbers are stored as text,
# matching the source data which keeps prices as formatted text strings (e.g. "69.426.93").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.426.93"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "3.770.79"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "614.90"
$ws.Range("E5").Value = "  -0.85%  "

$ws.Range("D6").Value = "176.91"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("D7").Value = "3.767.90"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  -1.46%  "

$ws.Range("D11").Value = "6.45"
$ws.Range("E11").Value = "  +2.66%  "

$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("D13").Value = "39.68"
$ws.Range("E13").Value = "  -4.07%  "

$ws.Range("D14").Value = "0.0000254"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").Value = "4.393.95"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "3.765.26"
$ws.Range("E16").Value = "  +0.25%  "

$ws.Range("D17").Value = "69.535.64"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "7.55"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("E19").Value = "  -3.56%  "

$ws.Range("D20").Value = "509.64"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "16.58"
$ws.Range("E21").Value = "  -1.06%  "

$ws.Range("D22").Value = "9.54"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("E24").Value = "  -1.62%  "

$ws.Range("D25").Value = "86.31"
$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D26").Value = "12.89"
$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("D27").Value = "0.0000141"
$ws.Range("E27").Value = "  +4.93%  "

$ws.Range("D28").Value = "10.53"
$ws.Range("E28").Value = "  -5.59%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").Value = "8.16"
$ws.Range("E32").Value = "  +2.69%  "

$ws.Range("D33").Value = "31.09"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  -0.22%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").Value = "6.12"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").Value = "0.142"
$ws.Range("E38").Value = "  +6.68%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "475.40"
$ws.Range("E39").Value = "  +11.06%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.340"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").Value = "2.06"
$ws.Range("E41").Value = "  -2.88%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "49.75"
$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.98"
$ws.Range("E43").Value = "  +5.46%  "

$ws.Range("D44").Value = "44.14"
$ws.Range("E44").Value = "  -3.39%  "

$ws.Range("D45").Value = "8.58"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").Value = "2.942.71"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("D47").Value = "0.0362"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").Value = "27.57"
$ws.Range("E48").Value = "  +0.15%  "

$ws.Range("D49").Value = "139.54"
$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").Value = "2.46"
$ws.Range("E51").Value = "  -2.24%  "
